# Apply updated crypto price/volume figures per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as literal text (matches the
# original inline-string cells) instead of auto-converting numeric-looking values
# like "314.83" or "1.000" into actual numbers.

$ws.Range("D2").Value = "'27.378.98"
$ws.Range("E2").Value = "'  +0.93%  "
$ws.Range("D3").Value = "'1.827.38"
$ws.Range("E3").Value = "'  +0.04%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'314.83"
$ws.Range("E5").Value = "'  +0.73%  "
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("D7").Value = "'0.4474"
$ws.Range("E7").Value = "'  -2.19%  "
$ws.Range("D8").Value = "'0.3792"
$ws.Range("E8").Value = "'  +1.68%  "
$ws.Range("D9").Value = "'0.07484"
$ws.Range("E9").Value = "'  +2.10%  "
$ws.Range("D10").Value = "'0.8879"
$ws.Range("E10").Value = "'  +3.22%  "
$ws.Range("E11").Value = "'  +0.30%  "
$ws.Range("D12").Value = "'1.831.21"
$ws.Range("E12").Value = "'  +0.32%  "
$ws.Range("D13").Value = "'6.762"
$ws.Range("E13").Value = "'  +0.99%  "
$ws.Range("D14").Value = "'5.460"
$ws.Range("E14").Value = "'  +2.17%  "
$ws.Range("D15").Value = "'93.87"
$ws.Range("E15").Value = "'  +1.03%  "
$ws.Range("D16").Value = "'0.07121"
$ws.Range("E16").Value = "'  +0.64%  "
$ws.Range("E17").Value = "'  -0.11%  "
$ws.Range("D18").Value = "'0.000008802"
$ws.Range("E18").Value = "'  -0.38%  "
$ws.Range("E19").Value = "'  -0.02%  "
$ws.Range("D20").Value = "'15.18"
$ws.Range("E20").Value = "'  +1.00%  "
$ws.Range("D21").Value = "'27.371.15"
$ws.Range("E21").Value = "'  +0.89%  "
$ws.Range("D22").Value = "'5.406"
$ws.Range("E22").Value = "'  +4.21%  "
$ws.Range("D23").Value = "'10.97"
$ws.Range("E23").Value = "'  -0.44%  "
$ws.Range("D24").Value = "'2.057.96"
$ws.Range("E24").Value = "'  +0.60%  "
$ws.Range("D25").Value = "'1.964"
$ws.Range("E25").Value = "'  -2.08%  "
$ws.Range("D26").Value = "'151.44"
$ws.Range("E26").Value = "'  -0.14%  "
$ws.Range("D27").Value = "'2.311"
$ws.Range("E27").Value = "'  +3.77%  "
$ws.Range("D28").Value = "'18.69"
$ws.Range("E28").Value = "'  +1.04%  "
$ws.Range("D29").Value = "'5.417"
$ws.Range("E29").Value = "'  +2.78%  "
$ws.Range("D30").Value = "'117.97"
$ws.Range("E30").Value = "'  +0.47%  "
$ws.Range("D31").Value = "'0.08893"
$ws.Range("E31").Value = "'  +0.29%  "
$ws.Range("D32").Value = "'0.7919"
$ws.Range("E32").Value = "'  +3.64%  "
$ws.Range("D33").Value = "'1.210"
$ws.Range("E33").Value = "'  +1.22%  "
$ws.Range("D34").Value = "'4.608"
$ws.Range("E34").Value = "'  +3.11%  "
$ws.Range("D35").Value = "'2.926"
$ws.Range("E35").Value = "'  -1.27%  "
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("D37").Value = "'1.111"
$ws.Range("E37").Value = "'  +0.48%  "
$ws.Range("D38").Value = "'0.01993"
$ws.Range("E38").Value = "'  +1.43%  "
$ws.Range("D39").Value = "'0.05319"
$ws.Range("E39").Value = "'  +0.56%  "
$ws.Range("D40").Value = "'7.335"
$ws.Range("E40").Value = "'  +2.23%  "
$ws.Range("D41").Value = "'0.5355"
$ws.Range("E41").Value = "'  -0.03%  "
$ws.Range("D42").Value = "'2.872"
$ws.Range("E42").Value = "'  -0.64%  "
$ws.Range("D43").Value = "'0.1724"
$ws.Range("E43").Value = "'  +0.46%  "
$ws.Range("D44").Value = "'2.328"
$ws.Range("E44").Value = "'  +16.23%  "
$ws.Range("D45").Value = "'8.687"
$ws.Range("E45").Value = "'  +0.76%  "
$ws.Range("D46").Value = "'0.5125"
$ws.Range("E46").Value = "'  -1.98%  "
$ws.Range("D47").Value = "'10.69"
$ws.Range("E47").Value = "'  -0.50%  "
$ws.Range("E48").Value = "'  +1.42%  "
$ws.Range("D49").Value = "'105.43"
$ws.Range("E49").Value = "'  -0.61%  "
$ws.Range("D50").Value = "'1.000"
$ws.Range("E50").Value = "'  +0.02%  "
$ws.Range("D51").Value = "'0.06408"
$ws.Range("E51").Value = "'  -1.17%  "
